$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 30, shifting existing rows 30-55 down to 31-56.
$ws.Rows(30).Insert()

# Populate the new row 30 with a new weekly price record (same dims/metadata
# as the former row 30, which is now row 31, but with updated date and prices).
$ws.Cells.Item(30, 1).Value = 1
$ws.Cells.Item(30, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(30, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(30, 4).Value = 44586
$ws.Cells.Item(30, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(30, 5).Value = 15
$ws.Cells.Item(30, 6).Value = 100112012
$ws.Cells.Item(30, 7).Value = "Espinaca"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 250
$ws.Cells.Item(30, 11).Value = 2500
$ws.Cells.Item(30, 12).Value = 3000
$ws.Cells.Item(30, 13).Value = 2750
$ws.Cells.Item(30, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item(30, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(30, 16).Value = 917
$ws.Cells.Item(30, 17).Value = 3
$ws.Cells.Item(30, 18).Value = "Hortaliza"
